# Add two new rows of test data (converted from an object[] in the source)
# to the invalidCredentialData worksheet, mirroring the existing rows'
# layout (Password/Language/ExpectedValue columns stay constant, only the
# "User Name" column value changes per row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "peter12"
$ws.Range("B4").Value = "peter123"
$ws.Range("C4").Value = "Danish"
$ws.Range("D4").Value = "Invalid username or password"

$ws.Range("A5").Value = "223frrr"
$ws.Range("B5").Value = "peter123"
$ws.Range("C5").Value = "Danish"
$ws.Range("D5").Value = "Invalid username or password"

# Leave the selection on the row following the newly-added data, matching
# where the cursor ended up after entering the rows in Excel.
$ws.Range("A6").Select()
